$d = $word.ActiveDocument

# The 4th table in the document is the "mcsa with n/a" results table
# (columns: N, Yes, Maybe, No, total_N, na - did not use / two rows:
# header + data). The edit adds a new left-hand label column holding
# the question text, with a blank header cell above it.
$t = $d.Tables.Item(4)

# Add a new first column to the table. This creates a new first cell
# in every existing row (row 1: header row, row 2: data row).
$t.Columns.Add() | Out-Null

# --- Row 1 (header row): the new first cell should just be an empty
# paragraph, matching the style of its neighboring header cells
# (bottom border + vertically bottom-aligned content).
$headerCell = $t.Cell(1, 1)
$headerCell.VerticalAlignment = 3      # wdCellAlignVerticalBottom
$headerCell.Borders.Item(-3).LineStyle = 1   # wdBorderBottom = single line

# --- Row 2 (data row): the new first cell should hold the question
# text, left aligned, using the "Compact" paragraph style (consistent
# with the rest of the table's body cells).
$dataCell = $t.Cell(2, 1)
$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:pPr><w:pStyle w:val='Compact'/><w:jc w:val='left'/></w:pPr>" +
       "<w:r><w:t xml:space='preserve'>Multiple choice with n/a option coded</w:t></w:r>" +
       "</w:p>"
$dataCell.Range.InsertXML($xml) | Out-Null
# InsertXML adds the new paragraph ahead of the cell's original
# (now-empty) paragraph rather than replacing it, so drop the leftover
# empty paragraph it pushed down.
$dataCell.Range.Paragraphs.Item(1).Range.Delete()

Write-Host "Added question-text label column to table 4."
